$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Hide the slide master shapes on this slide (adds showMasterSp="0" to p:sld)
$s.DisplayMasterShapes = 0

# Remove the stray "Rectangle 1" shape (id=2) that duplicated "Rectangle 14"
$s.Shapes.Item($s.Shapes.Count).Delete()
